# ReleaseNote_RT5514.docx - "New hotword model is merged." edit
#
# Summary of changes applied below (see commit message / diff):
#  1. Insert a brand-new release-note entry ("Version: 0.1.14") at the very
#     top of the document, above the existing "Version: 0.1.13" entry,
#     including its own Release date / Change / "1. New hotword model is
#     merged." paragraphs and a trailing blank paragraph. The "_GoBack"
#     bookmark moves from the old top paragraph into the new "... is
#     merged." paragraph (split right after "model i").
#  2. Remove the paragraph-mark run formatting (<w:pPr><w:rPr>...) from the
#     "1. Add boost control for both hw and sw path." paragraph.
#  3. Turn the following blank paragraph into a fully bare paragraph
#     (no pPr at all).
#  4. Add a <w:lastRenderedPageBreak/> to the start of the
#     "1. Merge "Farfield" hotword model into firmware." paragraph
#     (Version 0.0.8 Change section).
#  5. Remove the <w:lastRenderedPageBreak/> that used to sit at the start
#     of the "1.  Send an input event to user space when detect OKG."
#     paragraph (Version 0.0.7 Change section).

$d = $word.ActiveDocument
$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1. Insert the new "Version: 0.1.14" release-note block before the
#    current first paragraph ("Version: 0.1.13").
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs.First
$insertionPoint = $firstPara.Range.Duplicate
$insertionPoint.Collapse(1)  # wdCollapseStart

$newBlockXml = @"
<w:p $wdNS><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Version: 0.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="FF0000"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>.14</w:t></w:r></w:p>
<w:p $wdNS><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Release date</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>:2019</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>/5/24</w:t></w:r></w:p>
<w:p $wdNS><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Change:</w:t></w:r></w:p>
<w:p $wdNS><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">1. New </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hotword</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> model is merged.</w:t></w:r></w:p>
<w:p $wdNS><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>
"@

$insertionPoint.InsertXML($newBlockXml) | Out-Null

# Locate the freshly inserted "... is merged." paragraph (now paragraph 4)
# and split the bookmark into place right after "model i" / before
# "s merged." to match the target: "1. New hotword model i[_GoBack]s merged."
$mergedPara = $d.Paragraphs.Item(4)
$bmRange = $mergedPara.Range.Duplicate
$bmRange.Start = $mergedPara.Range.Start + "1. New hotword model i".Length
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Remove the bookmark from the old "Version: 0.1.13" paragraph (now shifted
# down, but still named/located via the unique bookmark name) -- the
# bookmark lived at the very end of that paragraph's text.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
# There are now two bookmarks named "_GoBack" (collection access by name
# returns the first / still-original one); delete it, leaving only the one
# we just created above.
if ($oldBookmark.Start -ne $bmRange.Start) {
    $oldBookmark.Delete()
} else {
    # Defensive fallback - if "Item" resolved to our new bookmark instead,
    # look the other one up by scanning all bookmarks.
    for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
        $cand = $d.Bookmarks.Item($i)
        if ($cand.Name -eq "_GoBack" -and $cand.Start -ne $bmRange.Start) {
            $cand.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# 2 & 3. "Add boost control" paragraph loses its paragraph-mark rFonts
#         hint, and the blank paragraph right after it becomes fully bare.
# ---------------------------------------------------------------------
$boostPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Add boost control for both*") {
        $boostPara = $p
        break
    }
}

$boostXml = @"
<w:p $wdNS><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">1. </w:t></w:r><w:r><w:t xml:space="preserve">Add boost control for both </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>hw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> path.</w:t></w:r></w:p>
"@
$boostPara.Range.InsertXML($boostXml) | Out-Null

# The blank paragraph right after it (index follows immediately).
$blankPara = $boostPara.Next()
$blankPara.Range.InsertXML("<w:p $wdNS/>") | Out-Null

# ---------------------------------------------------------------------
# 4. Add <w:lastRenderedPageBreak/> to the "Merge "Farfield" hotword
#    model into firmware." paragraph.
# ---------------------------------------------------------------------
$mergeFarfieldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*Merge*Farfield*hotword model into firmware*') {
        $mergeFarfieldPara = $p
        break
    }
}

$mergeFarfieldXml = @"
<w:p $wdNS><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">1. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Merge </w:t></w:r><w:r><w:t>&#x201C;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Farfield</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#x201D;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>hotword</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> model into firmware</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
"@
$mergeFarfieldPara.Range.InsertXML($mergeFarfieldXml) | Out-Null

# ---------------------------------------------------------------------
# 5. Remove <w:lastRenderedPageBreak/> from the "1.  Send an input event
#    to user space when detect OKG." paragraph.
# ---------------------------------------------------------------------
$okgPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*Send an input event to user space when*OKG*') {
        $okgPara = $p
        break
    }
}

$okgXml = @"
<w:p $wdNS><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">1.  </w:t></w:r><w:r><w:t xml:space="preserve">Send an input event to user space when </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>detect</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> OKG.</w:t></w:r></w:p>
"@
$okgPara.Range.InsertXML($okgXml) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
